$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are formatted as plain text (e.g. "40.138.39",
# "293.32") in the source sheet. Force text storage via NumberFormat "@"
# so Excel does not auto-coerce them into numeric cells, then restore the
# "Normal" style so no stray per-cell formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.138.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.234.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.22%  "

$ws.Range("E12").Value = "  +2.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.581.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("E15").Value = "  -3.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.220.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.733"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.068.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0889"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.24%  "

$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "

$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "155.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.66%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("E34").Value = "  -0.45%  "

$ws.Range("E35").Value = "  -1.97%  "

$ws.Range("E36").Value = "  +5.50%  "

$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0974"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.25%  "

$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.129.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.73%  "

$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.445.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "

# Row 51: coin entry replaced (BitcoinSV -> Aave)
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "88.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
